$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "61.005.78"
$ws.Range("E2").Value = "  +5.49%  "
$ws.Range("D3").Value = "2.380.30"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "549.89"
$ws.Range("E5").Value = "  +2.64%  "
Set-TextValue "D6" "134.67"
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +3.79%  "
$ws.Range("D9").Value = "2.379.18"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +2.95%  "
Set-TextValue "D14" "24.36"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "2.806.95"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "60.912.43"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "2.383.67"
$ws.Range("E18").Value = "  +4.36%  "
Set-TextValue "D19" "10.85"
Set-TextValue "D20" "4.25"
$ws.Range("E20").Value = "  +0.13%  "
Set-TextValue "D21" "6.95"
$ws.Range("E21").Value = "  +8.76%  "
Set-TextValue "D22" "319.53"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("E23").Value = "  +0.46%  "
Set-TextValue "D24" "63.77"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +4.87%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.07%  "
Set-TextValue "D27" "8.22"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  +4.92%  "
Set-TextValue "D31" "171.53"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +6.62%  "
Set-TextValue "D33" "5.98"
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("E34").Value = "  +14.61%  "
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("E37").Value = "  +0.02%  "
Set-TextValue "D38" "4.24"
$ws.Range("E38").Value = "  +8.78%  "
$ws.Range("E39").Value = "  +0.05%  "
Set-TextValue "D40" "328.95"
$ws.Range("E40").Value = "  +14.17%  "
Set-TextValue "D41" "1.58"
$ws.Range("E41").Value = "  +6.52%  "
Set-TextValue "D42" "38.58"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("E44").Value = "  +3.59%  "
Set-TextValue "D45" "0.0959"
$ws.Range("E45").Value = "  +1.53%  "
Set-TextValue "D46" "19.62"
$ws.Range("E46").Value = "  +7.76%  "
Set-TextValue "D47" "0.0505"
$ws.Range("E47").Value = "  +1.60%  "
Set-TextValue "D48" "0.569"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +5.07%  "
